# TeamOne-FinalProjectPlan.xlsx — "Updated Project Plan, added canvas to LaurenScene"
#
# Changes applied:
#   1. Insert a new task row (row 10) for "Find Player model, add all
#      animations" (Lauren, due 43775, Done), pushing the existing rows
#      10-22 down to 11-23.
#   2. Mark several tasks as completed / in progress:
#        - Set up Unity and Github projects         -> Done
#        - Find environment and theme assets        -> Done
#        - Find mobs and boss models with animations -> Done
#        - Program Character Controller (Movement)  -> Done
#        - Program Enemy AI states                  -> In Progress
#   3. Update the sheet's active selection to D12 (scrolled near row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Find Player model, add all animations" row -------
# Insert a blank row at position 10 (existing row 9's format/height is
# cloned onto it so borders / number formats / row height match the rest
# of the table).
$ws.Rows.Item(10).Insert()
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Rows.Item(10).RowHeight = 30

$ws.Range("A10").Value = "Find Player model, add all animations"
$ws.Range("B10").Value = 43775
$ws.Range("C10").Value = "Lauren"
$ws.Range("D10").Value = "Done"

# --- 2. Update task statuses ----------------------------------------------
$ws.Range("D7").Value = "Done"           # Set up Unity and Github projects
$ws.Range("D8").Value = "Done"           # Find environment and theme assets
$ws.Range("D9").Value = "Done"           # Find mobs and boss models with animations
$ws.Range("D13").Value = "Done"          # Program Character Controller (Movement)
$ws.Range("D16").Value = "In Progress"   # Program Enemy AI states

# --- 3. Update the view / selection ---------------------------------------
$ws.Range("D12").Select()
$excel.ActiveWindow.ScrollRow = 3

Write-Output "Applied Project Plan updates"
